$wb = $excel.ActiveWorkbook

# --- Sheet "Workspace": Business Line changed from "Corporate BE" to "COUNTERPARTY_BIB" ---
$wsWorkspace = $wb.Worksheets.Item("Workspace")
$wsWorkspace.Activate()
$wsWorkspace.Range("C3").Value = "COUNTERPARTY_BIB"
$wsWorkspace.Range("D3").Value = "COUNTERPARTY_BIB"
$wsWorkspace.Range("D3").Select() | Out-Null

# --- Sheet "r Workspace_AnalysisUnit": Id/Name changed from "WS000001_CUSTOMER_BE" to "WS000001_COUNTERPARTY_BIB" ---
$wsAnalysisUnit = $wb.Worksheets.Item("r Workspace_AnalysisUnit")
$wsAnalysisUnit.Activate()
$wsAnalysisUnit.Range("B3").Value = "WS000001_COUNTERPARTY_BIB"
$wsAnalysisUnit.Range("C3").Value = "WS000001_COUNTERPARTY_BIB"
$wsAnalysisUnit.Range("D3").Value = "WS000001_COUNTERPARTY_BIB"
$wsAnalysisUnit.Range("F3").Select() | Out-Null

# --- Sheet "r Workspace_TargetVariable" becomes the active tab ---
$wsTargetVariable = $wb.Worksheets.Item("r Workspace_TargetVariable")
$wsTargetVariable.Activate()
$wsTargetVariable.Range("J12").Select() | Out-Null
